$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing time-log entry on row 89
$ws.Range("D89").Value = 20
$ws.Range("A89").Value = 41930
$ws.Range("B89").Value = 0.63888888888888895
$ws.Range("C89").Value = 0.74652777777777779
$ws.Range("F89").Value = "Coding"

# Update selection to reflect the recorded edit position
$ws.Activate()
$ws.Range("H88").Select()
